$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 271
$targetDate = [DateTime]::FromOADate(45205)
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    if ($current -eq $targetDate -or $current -eq 45205) {
        $cell.Value = 45206
    }
}
